# Updates the crypto price/volume columns (D, E) per the latest scrape.
# Values that parse as plain numbers get a leading apostrophe so Excel
# keeps them as literal text (matching the source data, which stores
# every price/volume cell as a string).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.351.30'
$ws.Range("E2").Value = '  +0.48%  '
$ws.Range("D3").Value = '1.592.56'
$ws.Range("E3").Value = '  +0.61%  '
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").Value = '''211.78'
$ws.Range("E5").Value = '  +0.93%  '
$ws.Range("D6").Value = '''0.503'
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("D9").Value = '''0.0610'
$ws.Range("D10").Value = '''19.41'
$ws.Range("E10").Value = '  -0.72%  '
$ws.Range("E11").Value = '  +0.12%  '
$ws.Range("D12").Value = '1.816.60'
$ws.Range("E12").Value = '  +0.62%  '
$ws.Range("D13").Value = '1.624.28'
$ws.Range("E13").Value = '  +2.27%  '
$ws.Range("E14").Value = '  +0.61%  '
$ws.Range("D15").Value = '''0.522'
$ws.Range("E15").Value = '  +0.79%  '
$ws.Range("D16").Value = '''64.58'
$ws.Range("E16").Value = '  -0.18%  '
$ws.Range("D17").Value = '26.357.12'
$ws.Range("E17").Value = '  +0.52%  '
$ws.Range("E18").Value = '  -0.83%  '
$ws.Range("D19").Value = '''7.48'
$ws.Range("E19").Value = '  +3.50%  '
$ws.Range("D20").Value = '''212.73'
$ws.Range("E20").Value = '  +2.78%  '
$ws.Range("E21").Value = '  -0.26%  '
$ws.Range("E22").Value = '  +1.01%  '
$ws.Range("E23").Value = '  -1.23%  '
$ws.Range("D24").Value = '''9.00'
$ws.Range("E24").Value = '  +1.78%  '
$ws.Range("D25").Value = '''144.82'
$ws.Range("E25").Value = '  +0.13%  '
$ws.Range("E27").Value = '  +0.59%  '
$ws.Range("E28").Value = '  -0.67%  '
$ws.Range("D29").Value = '''15.20'
$ws.Range("E29").Value = '  -0.33%  '
$ws.Range("D30").Value = '''0.0502'
$ws.Range("E30").Value = '  -0.16%  '
$ws.Range("E31").Value = '  +0.84%  '
$ws.Range("D32").Value = '''3.21'
$ws.Range("E32").Value = '  -0.28%  '
$ws.Range("D33").Value = '''2.98'
$ws.Range("E33").Value = '  +1.49%  '
$ws.Range("D34").Value = '1.343.75'
$ws.Range("E34").Value = '  +4.33%  '
$ws.Range("E35").Value = '  -1.09%  '
$ws.Range("D36").Value = '''0.603'
$ws.Range("E36").Value = '  -0.21%  '
$ws.Range("E37").Value = '  +0.33%  '
$ws.Range("E38").Value = '  +0.06%  '
$ws.Range("E39").Value = '  -18.41%  '
$ws.Range("D40").Value = '''0.819'
$ws.Range("E40").Value = '  +0.54%  '
$ws.Range("D41").Value = '''5.78'
$ws.Range("E41").Value = '  +4.35%  '
$ws.Range("E42").Value = '  -0.24%  '
$ws.Range("E43").Value = '  +0.19%  '
$ws.Range("E44").Value = '  -0.86%  '
$ws.Range("D45").Value = '1.729.13'
$ws.Range("E45").Value = '  +0.57%  '
$ws.Range("D46").Value = '''61.63'
$ws.Range("E46").Value = '  -1.34%  '
$ws.Range("D47").Value = '''87.81'
$ws.Range("E47").Value = '  -1.12%  '
$ws.Range("E48").Value = '  +2.29%  '
$ws.Range("E49").Value = '  -2.89%  '
$ws.Range("D50").Value = '''0.0986'
$ws.Range("E50").Value = '  -3.24%  '
$ws.Range("D51").Value = '''0.0506'
$ws.Range("E51").Value = '  -0.60%  '
